$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (2-7) for the columns that changed, following the
# re-run of the NATMI pipeline with an updated expressing-cell count
# (1 -> 3) that cascades into the derived averages/totals/specificities.
$rows = @(
    @{ Row = 2;  E = 3; G = 0.3616933333333334; H = 1.08508;            K = 3; M = 0.6986340000000001;  N = 2.095902;             O = 0.0100360027308206;  P = 0.0100360027308206;  Q = 0.25269126024;     R = 2.27422134216;    S = 0.0100360027308206;  T = 0.0100360027308206 },
    @{ Row = 3;  E = 3; G = 0.3616933333333334; H = 1.08508;            K = 3; M = 45.91636366666668;    N = 137.749091;           O = 0.6595968005393648;  P = 0.6595968005393646;  Q = 16.60764262914223; R = 149.46878366228;   S = 0.6595968005393648;  T = 0.6595968005393646 },
    @{ Row = 4;  E = 3; G = 0.3616933333333334; H = 1.08508;            K = 3; M = 0.2330663333333333;   N = 0.699199;             O = 0.003348039685723394;P = 0.003348039685723394;Q = 0.08429853899111113; R = 0.7586868509200001; S = 0.003348039685723394;T = 0.003348039685723394 },
    @{ Row = 5;  E = 3; G = 0.3616933333333334; H = 1.08508;            K = 3; M = 0.4381680000000001;   N = 1.314504;             O = 0.006294361918484072;P = 0.006294361918484071;Q = 0.15848244448;       R = 1.42634200032;      S = 0.006294361918484072;T = 0.006294361918484071 },
    @{ Row = 6;  E = 3; G = 0.3616933333333334; H = 1.08508;            K = 3; M = 2.539916666666667;    N = 7.61975;              O = 0.03648635852638638; P = 0.03648635852638638; Q = 0.9186709255555556; R = 8.26803833;         S = 0.03648635852638638; T = 0.03648635852638638 },
    @{ Row = 7;  E = 3; G = 0.3616933333333334; H = 1.08508;            K = 3; M = 19.78662633333333;    N = 59.35987900000001;    O = 0.284238436599221;   P = 0.2842384365992209;  Q = 7.156690833924445;  R = 64.41021750532001;  S = 0.284238436599221;   T = 0.2842384365992209 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("E$row").Value = $r.E
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("K$row").Value = $r.K
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = $r.T
}
